# Add data for 2022-08-12:
# Roll the "through" date in the workbook from 2022-08-03 to 2022-08-04,
# and update the August (row 9) and Total (row 14) figures in the
# "2022 (through ...)" column (I) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new as-of date.
$ws.Name = "Through 2022-08-04"

# Update the "Total" column header text (shared string used by I1).
$ws.Range("I1").Value = "2022 (through 08-04)"

# August 2022 year-to-date count increases from 16 to 22.
$ws.Range("I9").Value = 22

# Grand total increases from 986 to 992.
$ws.Range("I14").Value = 992
